$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rapor Sayfası")
$ws2 = $wb.Worksheets.Item("M_Listesi")

# Date/time value
$ws1.Range("N1").Value = 44337.86379885743

# Power values
$ws1.Range("B3").Value = "18774.0 VA"
$ws1.Range("H3").Value = "1.772 kg"
$ws1.Range("H7").Value = "Al:1.772 //  Cu:0.0kg"

$ws1.Range("G10").Value = "447.0 V ( 447.0V )"
$ws1.Range("K10").Value = "Boşluk"

$ws1.Range("L11").Value = 56.462
$ws1.Range("L14").Value = 2900.122
$ws1.Range("L15").Value = 45
$ws1.Range("L16").Value = 50.03
$ws1.Range("B17").Value = 44.55411000000001
$ws1.Range("L18").Value = 107.639
$ws1.Range("G19").Value = "129 sp"
$ws1.Range("L19").Value = 0.129
$ws1.Range("B20").Value = 0.5905020852980082
$ws1.Range("L20").Value = 1.14
$ws1.Range("L21").Value = 100.792
$ws1.Range("L22").Value = 0.02
$ws1.Range("L23").Value = 1.023
$ws1.Range("L24").Value = 37.999
$ws1.Range("L27").Value = 128.724

$ws1.Range("N11").Select()

$ws2.Range("D3").Value = 0.5905020852980082
